$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("testT2922", iteration 2) is deleted entirely; rows 11-20 shift up to 10-19.
$ws.Rows.Item(10).Delete()

# A new row 20 is appended ("testT4244", iteration 1), formatted like the other
# "iteration 1 / YES" rows (e.g. row 17, which carries style index 3: font2/no-fill).
$ws.Range("A17:C17").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A20").Value = "testT4244"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "YES"

# Update the saved selection/active cell.
$ws.Range("A7").Select()
